{"js": "// Replace the text of each lattice-multiplication exercise cell, in\n// document order, with the new problem data while keeping the existing\n// paragraph/run/formatting (sz=32) and the manual line breaks between\n// the 5 lines inside every cell.\n//\n// Word represents a manual line break (<w:br/>) as \"\\u000b\" (vertical\n// tab) in Office.js Range.text / insertText, so joining the 5 lines of\n// each cell with \"\\u000b\" reproduces the <w:t>...</w:t><w:br/> pattern\n// exactly.\n\nconst NEW_CELLS = [\n  [\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"],\n  [\"87 x 35\", \"  3    5\", \"  ----\", \"8|    |\", \"7|    |\"],\n  [\"51 x 86\", \"  8    6\", \"  ----\", \"5|    |\", \"1|    |\"],\n  [\"20 x 82\", \"  8    2\", \"  ----\", \"2|    |\", \"0|    |\"],\n  [\"89 x 45\", \"  4    5\", \"  ----\", \"8|    |\", \"9|    |\"],\n  [\"24 x 66\", \"  6    6\", \"  ----\", \"2|    |\", \"4|    |\"],\n  [\"93 x 74\", \"  7    4\", \"  ----\", \"9|    |\", \"3|    |\"],\n  [\"38 x 82\", \"  8    2\", \"  ----\", \"3|    |\", \"8|    |\"],\n  [\"56 x 93\", \"  9    3\", \"  ----\", \"5|    |\", \"6|    |\"],\n  [\"87 x 49\", \"  4    9\", \"  ----\", \"8|    |\", \"7|    |\"],\n  [\"11 x 89\", \"  8    9\", \"  ----\", \"1|    |\", \"1|    |\"],\n  [\"38 x 53\", \"  5    3\", \"  ----\", \"3|    |\", \"8|    |\"],\n  [\"59 x 33\", \"  3    3\", \"  ----\", \"5|    |\", \"9|    |\"],\n  [\"48 x 56\", \"  5    6\", \"  ----\", \"4|    |\", \"8|    |\"],\n  [\"31 x 39\", \"  3    9\", \"  ----\", \"3|    |\", \"1|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    if (i >= NEW_CELLS.length) break;\n    const lines = NEW_CELLS[i];\n    const text = lines.join(\"\\u000b\");\n    cell.getRange().insertText(text, \"Replace\");\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the text of each lattice-multiplication exercise cell, in\n# document order, with the new problem data while keeping the existing\n# paragraph/run/formatting (sz=32) and the manual line breaks between\n# the 5 lines inside every cell.\n#\n# Word represents a manual line break (<w:br/>) as Chr(11) (vertical\n# tab) in Range.Text, so joining the 5 lines of each cell with Chr(11)\n# reproduces the <w:t>...</w:t><w:br/> pattern exactly. Assigning to\n# Range.Text (rather than InsertAfter) replaces the cell's whole\n# paragraph content but leaves the end-of-cell mark untouched.\n\n$vt = [char]11\n\n$newCells = @(\n    @(\"75 x 76\", \"  7    6\", \"  ----\", \"7|    |\", \"5|    |\"),\n    @(\"87 x 35\", \"  3    5\", \"  ----\", \"8|    |\", \"7|    |\"),\n    @(\"51 x 86\", \"  8    6\", \"  ----\", \"5|    |\", \"1|    |\"),\n    @(\"20 x 82\", \"  8    2\", \"  ----\", \"2|    |\", \"0|    |\"),\n    @(\"89 x 45\", \"  4    5\", \"  ----\", \"8|    |\", \"9|    |\"),\n    @(\"24 x 66\", \"  6    6\", \"  ----\", \"2|    |\", \"4|    |\"),\n    @(\"93 x 74\", \"  7    4\", \"  ----\", \"9|    |\", \"3|    |\"),\n    @(\"38 x 82\", \"  8    2\", \"  ----\", \"3|    |\", \"8|    |\"),\n    @(\"56 x 93\", \"  9    3\", \"  ----\", \"5|    |\", \"6|    |\"),\n    @(\"87 x 49\", \"  4    9\", \"  ----\", \"8|    |\", \"7|    |\"),\n    @(\"11 x 89\", \"  8    9\", \"  ----\", \"1|    |\", \"1|    |\"),\n    @(\"38 x 53\", \"  5    3\", \"  ----\", \"3|    |\", \"8|    |\"),\n    @(\"59 x 33\", \"  3    3\", \"  ----\", \"5|    |\", \"9|    |\"),\n    @(\"48 x 56\", \"  5    6\", \"  ----\", \"4|    |\", \"8|    |\"),\n    @(\"31 x 39\", \"  3    9\", \"  ----\", \"3|    |\", \"1|    |\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($i -ge $newCells.Count) { continue }\n        $lines = $newCells[$i]\n        $text = [string]::Join($vt, $lines)\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $text\n        $i++\n    }\n}\n"}
